$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting of the other
# header cells (B1:G1) by copying G1's formatting over to H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the data values for the new column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
